$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values
# Column D values must stay as text (matching original inlineStr cells), so we
# temporarily force a Text number format before assignment and then restore the
# default "Normal" style so no extra style index is left on the cell.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "56.548.03"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -3.19%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.978.05"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -5.63%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "496.26"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -5.79%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.11"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  -0.14%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.425"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -3.45%  "
$ws.Range("E11").Value = "  -7.63%  "
$ws.Range("E12").Value = "  -0.69%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.484.19"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -5.81%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.14"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.26%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "56.469.96"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.33%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.975.99"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("E17").Value = "  -3.83%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "5.81"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.33"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.60%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.79"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "325.83"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("E22").Value = "  -0.35%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.469"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -8.13%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "61.57"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -8.05%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -5.99%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0₃0894"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.50%  "
$ws.Range("E28").Value = "  -0.05%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.24%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.76"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.50%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.74"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.87%  "
$ws.Range("E32").Value = "  -7.14%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "20.28"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -5.50%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "152.33"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.39%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.45"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -8.32%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.28"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -7.43%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.60"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -10.70%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0670"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.53%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "23.14"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.73%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.002.71"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.82%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.62"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -9.56%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.640"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -8.16%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -9.17%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.221.72"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.87%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.40"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -9.35%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.44%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0236"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.68%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -7.06%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "18.98"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -8.58%  "
